$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update rotates the price-detail columns (D, L, M, N, O, P, Q, S, T)
# among rows 2-7, while the descriptive columns stay put.
# Mapping: new row -> source row (old values to copy in)
# 2 <- 4, 3 <- 2, 4 <- 5, 5 <- 3, 6 <- 7, 7 <- 6

$source = @{
    2 = @(44392, 'Especial',     500, 7000,  7000,  7000,  '$/bandeja 8 kilos',        875,  8)
    3 = @(44208, 'Especial',      70, 24000, 24000, 24000, '$/caja 15 kilos granel',   1600, 15)
    4 = @(44217, 'Primera',       55, 18000, 18000, 18000, '$/caja 18 kilos granel',   1000, 18)
    5 = @(44418, 'Especial',     100, 8000,  8000,  8000,  '$/caja 15 kilos granel',    533, 15)
    6 = @(44427, 'Primera',       55, 7000,  7000,  7000,  '$/caja 15 kilos granel',    467, 15)
    7 = @(44264, 'Calibre 100',   50, 20000, 20000, 20000, '$/caja 18 kilos embalada', 1111, 18)
}

foreach ($row in $source.Keys) {
    $vals = $source[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D Fecha
    $ws.Cells.Item($row, 12).Value = $vals[1]   # L Calidad
    $ws.Cells.Item($row, 13).Value = $vals[2]   # M Volumen
    $ws.Cells.Item($row, 14).Value = $vals[3]   # N Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[4]   # O Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals[6]   # Q Unidad de comercializacion
    $ws.Cells.Item($row, 19).Value = $vals[7]   # S Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals[8]   # T Kg / unidad
}
